$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows right before the current row 967 (pushes existing
# rows 967-1015 down to 970-1018), matching the new dimension A1:T1018.
$ws.Range("A967:A969").EntireRow.Insert()

# Columns A, B, C, E-K are identical for every record in this block
# (Mercado ID, Mercado, Región, Codreg, Tipo, Producto ID, Producto,
# Categoría ID, Categoría, Variedad). Populate the three new rows explicitly
# (row by row, since COM PasteSpecial tiling is unreliable right after an
# Insert in this environment).
foreach ($r in 967..969) {
    $ws.Range("A$r").Value = 9
    $ws.Range("B$r").Value = "Vega Central Mapocho de Santiago"
    $ws.Range("C$r").Value = "Metropolitana"
    $ws.Range("E$r").Value = 13
    $ws.Range("F$r").Value = "Fruta"
    $ws.Range("G$r").Value = 100101
    $ws.Range("H$r").Value = "Berries"
    $ws.Range("I$r").Value = 100101007
    $ws.Range("J$r").Value = "Kiwi"
    $ws.Range("K$r").Value = "Hayward"
}

# Row 967 - Especial
$ws.Range("D967").Value = 45147
$ws.Range("L967").Value = "Especial"
$ws.Range("M967").Value = 220
$ws.Range("N967").Value = 8000
$ws.Range("O967").Value = 8000
$ws.Range("P967").Value = 8000
$ws.Range("Q967").Value = "$/bandeja 10 kilos"
$ws.Range("R967").Value = "Región de O'Higgins"
$ws.Range("S967").Value = 800
$ws.Range("T967").Value = 10

# Row 968 - Primera
$ws.Range("D968").Value = 45147
$ws.Range("L968").Value = "Primera"
$ws.Range("M968").Value = 350
$ws.Range("N968").Value = 6000
$ws.Range("O968").Value = 6000
$ws.Range("P968").Value = 6000
$ws.Range("Q968").Value = "$/bandeja 10 kilos"
$ws.Range("R968").Value = "Región de O'Higgins"
$ws.Range("S968").Value = 600
$ws.Range("T968").Value = 10

# Row 969 - Segunda
$ws.Range("D969").Value = 45147
$ws.Range("L969").Value = "Segunda"
$ws.Range("M969").Value = 220
$ws.Range("N969").Value = 4000
$ws.Range("O969").Value = 4000
$ws.Range("P969").Value = 4000
$ws.Range("Q969").Value = "$/bandeja 10 kilos"
$ws.Range("R969").Value = "Región de O'Higgins"
$ws.Range("S969").Value = 400
$ws.Range("T969").Value = 10
